# Fix column name and add drop down option
# - Rename the "Clinical Notes" header (column O) to "Clinic Notes"
# - Reset the sheet's active selection back to the (now visible) O1 cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mis-spelled / mis-named column header in O1.
$ws.Range("O1").Value = "Clinic Notes"

# Bring the view back to show column O / select O1 (undoing the prior
# scrolled-right view that had topLeftCell=U1 and selection X2).
$ws.Range("O1").Select() | Out-Null
